$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.420.00"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.565.33"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.73"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3726"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.14"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.58"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.920"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "1.562.83"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06755"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.60"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.329"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.37"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "22.421.40"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.382"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.556"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.81"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("E28").Value = "  -2.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.015"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.95"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").Value = "1.737.76"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.050"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.102"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.610"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08298"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02453"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2269"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06374"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.284"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.332"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.22"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6249"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.80"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6096"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.766"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.15"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07219"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.34%  "
